$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 4 (old rows 4-7 shift down to 5-8),
# carrying over the row-4 formatting (date style) as Excel normally does.
$ws.Rows("4:4").Insert()

# Populate the newly inserted row 4 with the new weekly data point.
$ws.Range("A4").Value = 7
$ws.Range("B4").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C4").Value = "Ñuble"
$ws.Range("D4").Value = 45028
$ws.Range("E4").Value = 16
$ws.Range("F4").Value = "Fruta"
$ws.Range("G4").Value = 100107
$ws.Range("H4").Value = "Otros"
$ws.Range("I4").Value = 100107011
$ws.Range("J4").Value = "Tuna"
$ws.Range("K4").Value = "Sin especificar"
$ws.Range("L4").Value = "Primera"
$ws.Range("M4").Value = 50
$ws.Range("N4").Value = 18000
$ws.Range("O4").Value = 18000
$ws.Range("P4").Value = 18000
$ws.Range("Q4").Value = "$/caja 18 kilos"
$ws.Range("R4").Value = "Región Metropolitana"
$ws.Range("S4").Value = 1000
$ws.Range("T4").Value = 18
